$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The extractor now normalizes to a fixed-width result set, so the table
# gains one more data row (Row=4). Append it directly below the existing
# data, matching the layout/format already used by the row above it.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 34234
$ws.Range("C5").Value = "jj100 device"
$ws.Range("D5").Value = "JM100"
$ws.Range("E5").Value = "JJ150"
$ws.Range("F5").Value = 41092
$ws.Range("F5").NumberFormat = "MM/DD/YY"
